{"js": "// Remove the trailing \"Ver no Jupiter Salvar em pdf Salvar em docx\" paragraph,\n// the \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github\n// pages. Original theme under Creative Commons Attribution\" paragraph, and\n// the blank paragraph that separates them from the preceding requirement\n// line (\"LOB1038: F\u00edsica Experimental I (Requisito fraco)\").\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the paragraph that immediately precedes the block to be removed.\nlet anchorIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.indexOf(\"LOB1038: F\u00edsica Experimental I (Requisito fraco)\") !== -1) {\n    anchorIndex = i;\n    break;\n  }\n}\n\nif (anchorIndex !== -1 && anchorIndex + 3 <= items.length - 1) {\n  const blank = items[anchorIndex + 1];\n  const jupiter = items[anchorIndex + 2];\n  const copyright = items[anchorIndex + 3];\n\n  if (\n    blank.text.trim() === \"\" &&\n    jupiter.text.indexOf(\"Ver no Jupiter Salvar em pdf Salvar em docx\") !== -1 &&\n    copyright.text.indexOf(\"Powered by Jekyll and Github pages\") !== -1\n  ) {\n    // Delete in reverse order so earlier deletions never invalidate the\n    // object references of the paragraphs still queued for removal.\n    copyright.delete();\n    jupiter.delete();\n    blank.delete();\n\n    await context.sync();\n  }\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the \"LOB1038: F\u00edsica Experimental I (Requisito fraco)\" paragraph -\n# the blank paragraph, the \"Ver no Jupiter...\" paragraph and the\n# \"\u00a9 2020 ...\" footer paragraph that immediately follow it must be removed.\n$anchorIndex = -1\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($t -like \"*LOB1038: F\u00edsica Experimental I (Requisito fraco)*\") {\n        $anchorIndex = $i\n        break\n    }\n}\n\nif ($anchorIndex -ne -1 -and ($anchorIndex + 3) -le $count) {\n    $blankText = $d.Paragraphs.Item($anchorIndex + 1).Range.Text.Trim()\n    $jupiterText = $d.Paragraphs.Item($anchorIndex + 2).Range.Text\n    $copyrightText = $d.Paragraphs.Item($anchorIndex + 3).Range.Text\n\n    if ($blankText -eq \"\" -and\n        $jupiterText -like \"*Ver no Jupiter Salvar em pdf Salvar em docx*\" -and\n        $copyrightText -like \"*Powered by Jekyll and Github pages*\") {\n\n        # Delete from the bottom up so the indices of paragraphs still queued\n        # for removal remain valid while earlier ones disappear.\n        $d.Paragraphs.Item($anchorIndex + 3).Range.Delete()\n        $d.Paragraphs.Item($anchorIndex + 2).Range.Delete()\n        $d.Paragraphs.Item($anchorIndex + 1).Range.Delete()\n    }\n}\n"}
